# asignamos WoS en base a clase social
# Add three new rows (Salariat, Intermediate, Working) with the same
# schedule values as the existing "work" row (row 4), and leave an
# extra blank row (14) below them, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNames = @("Salariat", "Intermediate", "Working")

$destRow = 11
foreach ($name in $newNames) {
    # Copy the full "work" row (row 4) formatting + values into the new row
    $srcRange = $ws.Range("A4:R4")
    $dstRange = $ws.Range("A" + $destRow + ":R" + $destRow)
    $srcRange.Copy($dstRange)

    # Overwrite the archetype name in column A with the new social-class name
    $ws.Cells.Item($destRow, 1).Value = $name

    $destRow++
}

# Trailing (empty) row, keeping the same column-A formatting as the rest
$ws.Range("A4").Copy($ws.Cells.Item($destRow, 1))
$ws.Cells.Item($destRow, 1).ClearContents()

# Update selection/pane to match the final saved view
[void]$ws.Range("C13:R13").Select()
